# Add updated pipe fluid flow figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15.23

$ws.Range("A2").Value = 28.630468503753782
$ws.Range("B2").Value = 27.933647815573622

$ws.Range("A3").Value = 46.424324163695431
$ws.Range("B3").Value = 34.955264023511084

$ws.Range("A4").Value = 59.220651167493827
$ws.Range("B4").Value = 39.099598143212624

$ws.Range("A5").Value = 73.521852300752712
$ws.Range("B5").Value = 43.079867294704265

$ws.Range("A6").Value = 90.011373845611956
$ws.Range("B6").Value = 46.935324028156764

$ws.Range("A7").Value = 109.33485710037152
$ws.Range("B7").Value = 50.610054172001263

$ws.Range("A8").Value = 132.21472731702607
$ws.Range("B8").Value = 53.978288755798374

$ws.Range("A9").Value = 159.5141169572409
$ws.Range("B9").Value = 56.838824487449905

$ws.Range("A10").Value = 192.29512052321328
$ws.Range("B10").Value = 58.894973353295178

$ws.Range("A11").Value = 231.88599055908699
$ws.Range("B11").Value = 59.721773696695422

$ws.Range("A1:B11").Select()
